# Weekly update: insert a new price-record row for
# "Feria Lagunitas de Puerto Montt - Pepino ensalada" at row 387,
# pushing the previously existing rows 387-406 down to 388-407.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 387 (shifts rows 387:406 -> 388:407)
$ws.Rows.Item(387).Insert()

# Populate the new row 387 with the new weekly record
$ws.Cells.Item(387, 1).Value = 4
$ws.Cells.Item(387, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(387, 3).Value = "Los Lagos"
$ws.Cells.Item(387, 4).Value = 45008
$ws.Cells.Item(387, 5).Value = 10
$ws.Cells.Item(387, 6).Value = 100112043
$ws.Cells.Item(387, 7).Value = "Pepino ensalada"
$ws.Cells.Item(387, 8).Value = "Sin especificar"
$ws.Cells.Item(387, 9).Value = "Primera"
$ws.Cells.Item(387, 10).Value = 200
$ws.Cells.Item(387, 11).Value = 13000
$ws.Cells.Item(387, 12).Value = 15000
$ws.Cells.Item(387, 13).Value = 14000
$ws.Cells.Item(387, 14).Value = "`$/caja 60 unidades"
$ws.Cells.Item(387, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(387, 16).Value = 233
$ws.Cells.Item(387, 17).Value = 60
$ws.Cells.Item(387, 18).Value = "Hortaliza"
